$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.768.88"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.305.63"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.55"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.98"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").Value = "2.304.01"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.75"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "2.717.18"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "58.638.42"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "2.300.36"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.16"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.17"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.62"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.42"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.92"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.30"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.97"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("D31").Value = "0.0₃0733"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.86"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.387"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("E37").Value = "  +4.23%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "296.11"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.25"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0959"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.64"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  +0.59%  "
